$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Docente(s)*") {
        $targetIndex = $i
    }
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new paragraph right after it.
$target.Range.InsertParagraphAfter()

# The freshly created paragraph is now the next one; fill it in and
# restyle it as a bulleted list item holding the professor's name.
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "5817650 - " + [char]0x00C9 + "rica Leonor Rom" + [char]0x00E3 + "o"
$newPara.Style = "List Bullet"
